$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Replace("Aida's Astounding Team", "A's Astounding Team")
$ws.Cells.Replace("Kelly's Deluxe Team", "K's Deluxe Team")
$ws.Cells.Replace("Magic Mikaela", "Magic M")
